# Add a new slide at the end of the deck.
#
# The source deck has 10 slides; the author appended an 11th slide built
# from the "Titulo e conteudo" (Title and Content) layout -- the same
# layout used elsewhere in presentation.xml as CustomLayouts item 2 -- and
# left its title / content placeholders empty.
$p = $ppt.ActivePresentation

$newIndex  = $p.Slides.Count + 1
$ppLayoutText = 2

$slide = $p.Slides.Add($newIndex, $ppLayoutText)

# Leave the title and content placeholders untouched (blank), matching the
# freshly inserted, still-empty slide from the authored edit.
